# Rebuild the workbook's "generator" columns being dropped, and refresh every
# numeric value to the new build's recomputed results (0294d82).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) info_solution: comp_time refreshed
# ---------------------------------------------------------------------------
$wsInfo = $wb.Worksheets.Item("info_solution")
$wsInfo.Range("A2").Value = 3.6772587299346924

# ---------------------------------------------------------------------------
# 2) economics_aggregator: row of aggregator economics refreshed
# ---------------------------------------------------------------------------
$wsEconAgg = $wb.Worksheets.Item("economics_aggregator")
$wsEconAgg.Range("B2").Value = -50703.11979608213
$wsEconAgg.Range("C2").Value = 45181.850661596065
$wsEconAgg.Range("D2").Value = -1206493.1898647628
$wsEconAgg.Range("E2").Value = -1206493.1898647628
$wsEconAgg.Range("F2").Value = 3036.930062754793

# ---------------------------------------------------------------------------
# 3) peak_aggregator: row of aggregator peak values refreshed
# ---------------------------------------------------------------------------
$wsPeakAgg = $wb.Worksheets.Item("peak_aggregator")
$wsPeakAgg.Range("B2").Value = 70.55334114690362
$wsPeakAgg.Range("C2").Value = 97.89381112765514
$wsPeakAgg.Range("D2").Value = 78.28045053530639
$wsPeakAgg.Range("E2").Value = 68.20875428947097
$wsPeakAgg.Range("F2").Value = 58.7148770548664
$wsPeakAgg.Range("H2").Value = 78.46006385188183
$wsPeakAgg.Range("I2").Value = 66.42216904631434
$wsPeakAgg.Range("K2").Value = 61.713700143020986
$wsPeakAgg.Range("M2").Value = 86.32036308984883

# ---------------------------------------------------------------------------
# 4) design_users: the "x_us_generator" column (F) is dropped entirely; the
#    remaining design columns shift left, and the surviving numbers refresh.
# ---------------------------------------------------------------------------
$wsDesign = $wb.Worksheets.Item("design_users")
$wsDesign.Range("F1:F4").EntireColumn.Delete()

# Refresh the values now sitting in F (was x_us_batt / G) and E (x_us_PV)
$wsDesign.Range("E4").Value = 22.68494581951114
$wsDesign.Range("F3").Value = 3.3348083265571247
$wsDesign.Range("G3").Value = 3.3348083265571247
$wsDesign.Range("F4").Value = 5.0986036118951
$wsDesign.Range("G4").Value = 5.0986036118951
$wsDesign.Range("H4").Value = 73.4786775552479

# ---------------------------------------------------------------------------
# 5) economics_users: "CAPEX_us_generator" (M) and "C_OEM_us_generator" (S)
#    columns are dropped entirely; everything to their right shifts left.
# ---------------------------------------------------------------------------
$wsEconUsers = $wb.Worksheets.Item("economics_users")
$wsEconUsers.Range("M1:M4").EntireColumn.Delete()
$wsEconUsers.Range("R1:R4").EntireColumn.Delete()   # was S before the M delete shifted it to R

# Refresh every surviving number to the new build's recomputed results
$wsEconUsers.Range("E2").Value = 0.0

$wsEconUsers.Range("B3").Value = -229654.1601644049
$wsEconUsers.Range("D3").Value = -10715.15628544633
$wsEconUsers.Range("G3").Value = 1352.476743299115
$wsEconUsers.Range("H3").Value = 7727.7803158796
$wsEconUsers.Range("I3").Value = 9274.367342711164
$wsEconUsers.Range("J3").Value = -138236.6653041894
$wsEconUsers.Range("M3").Value = 1333.9233306228498
$wsEconUsers.Range("N3").Value = 666.9616653114249
$wsEconUsers.Range("Q3").Value = 1400.1373529091506
$wsEconUsers.Range("R3").Value = 16.674041632785624
$wsEconUsers.Range("S3").Value = 6.669616653114249

$wsEconUsers.Range("B4").Value = -451303.9321147784
$wsEconUsers.Range("C4").Value = 259791.10814409857
$wsEconUsers.Range("D4").Value = -12566.864984448177
$wsEconUsers.Range("F4").Value = 43451.137106864444
$wsEconUsers.Range("G4").Value = 2067.8078417502916
$wsEconUsers.Range("I4").Value = 7253.824536721587
$wsEconUsers.Range("J4").Value = -143512.08077400186
$wsEconUsers.Range("L4").Value = 36295.913311217824
$wsEconUsers.Range("M4").Value = 2039.4414447580398
$wsEconUsers.Range("N4").Value = 1019.7207223790199
$wsEconUsers.Range("O4").Value = 220436.0326657437
$wsEconUsers.Range("Q4").Value = 680.5483745853342
$wsEconUsers.Range("R4").Value = 25.4930180594755
$wsEconUsers.Range("S4").Value = 10.1972072237902
$wsEconUsers.Range("T4").Value = 2204.360326657437

# ---------------------------------------------------------------------------
# 6) peak_users: minor refresh of the last two values
# ---------------------------------------------------------------------------
$wsPeakUsers = $wb.Worksheets.Item("peak_users")
$wsPeakUsers.Range("N2").Value = 22.572662993442876
$wsPeakUsers.Range("O2").Value = 21.737480146405954
